$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Log hours worked: Mon = 5 hours, Tues = 1.5 hours
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 1.5

# Update the selected/active cell to B4 (where next entry would go)
$ws.Range("B4").Select()
